$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Il34"
$ws.Range("C2").Value = "Csf1r"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 1.393664333333333
$ws.Range("H2").Value = 4.180993
$ws.Range("I2").Value = 0.04358216200908328
$ws.Range("J2").Value = 0.04358216200908328
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 22.07422366666666
$ws.Range("N2").Value = 66.22267099999999
$ws.Range("O2").Value = 0.8730897844203874
$ws.Range("P2").Value = 0.8730897844203874
$ws.Range("Q2").Value = 30.76405821025589
$ws.Range("R2").Value = 276.8765238923029
$ws.Range("S2").Value = 0.03805114043308493
$ws.Range("T2").Value = 0.03805114043308492

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Il34"
$ws.Range("C3").Value = "Csf1r"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 1.393664333333333
$ws.Range("H3").Value = 4.180993
$ws.Range("I3").Value = 0.04358216200908328
$ws.Range("J3").Value = 0.04358216200908328
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.456833333333333
$ws.Range("N3").Value = 4.3705
$ws.Range("O3").Value = 0.05762133790721463
$ws.Range("P3").Value = 0.05762133790721465
$ws.Range("Q3").Value = 2.030336656277778
$ws.Range("R3").Value = 18.2730299065
$ws.Range("S3").Value = 0.00251126248385236
$ws.Range("T3").Value = 0.00251126248385236

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Il34"
$ws.Range("C4").Value = "Csf1r"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 1.393664333333333
$ws.Range("H4").Value = 4.180993
$ws.Range("I4").Value = 0.04358216200908328
$ws.Range("J4").Value = 0.04358216200908328
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.751822333333333
$ws.Range("N4").Value = 5.255467
$ws.Range("O4").Value = 0.06928887767239802
$ws.Range("P4").Value = 0.06928887767239804
$ws.Range("Q4").Value = 2.441452304303445
$ws.Range("R4").Value = 21.973070738731
$ws.Range("S4").Value = 0.003019759092146004
$ws.Range("T4").Value = 0.003019759092146004

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Il34"
$ws.Range("C5").Value = "Csf1r"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 2.952820666666666
$ws.Range("H5").Value = 8.858462
$ws.Range("I5").Value = 0.092339529397755
$ws.Range("J5").Value = 0.092339529397755
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 22.07422366666666
$ws.Range("N5").Value = 66.22267099999999
$ws.Range("O5").Value = 0.8730897844203874
$ws.Range("P5").Value = 0.8730897844203874
$ws.Range("Q5").Value = 65.18122384355577
$ws.Range("R5").Value = 586.6310145920019
$ws.Range("S5").Value = 0.08062069981536595
$ws.Range("T5").Value = 0.08062069981536595

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Il34"
$ws.Range("C6").Value = "Csf1r"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 2.952820666666666
$ws.Range("H6").Value = 8.858462
$ws.Range("I6").Value = 0.092339529397755
$ws.Range("J6").Value = 0.092339529397755
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.456833333333333
$ws.Range("N6").Value = 4.3705
$ws.Range("O6").Value = 0.05762133790721463
$ws.Range("P6").Value = 0.05762133790721465
$ws.Range("Q6").Value = 4.301767574555555
$ws.Range("R6").Value = 38.715908171
$ws.Range("S6").Value = 0.00532072722562122
$ws.Range("T6").Value = 0.005320727225621221

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Il34"
$ws.Range("C7").Value = "Csf1r"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 2.952820666666666
$ws.Range("H7").Value = 8.858462
$ws.Range("I7").Value = 0.092339529397755
$ws.Range("J7").Value = 0.092339529397755
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.751822333333333
$ws.Range("N7").Value = 5.255467
$ws.Range("O7").Value = 0.06928887767239802
$ws.Range("P7").Value = 0.06928887767239804
$ws.Range("Q7").Value = 5.172817190194889
$ws.Range("R7").Value = 46.555354711754
$ws.Range("S7").Value = 0.006398102356767847
$ws.Range("T7").Value = 0.006398102356767849

$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Il34"
$ws.Range("C8").Value = "Csf1r"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 27.63137633333334
$ws.Range("H8").Value = 82.894129
$ws.Range("I8").Value = 0.8640783085931617
$ws.Range("J8").Value = 0.8640783085931617
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 22.07422366666666
$ws.Range("N8").Value = 66.22267099999999
$ws.Range("O8").Value = 0.8730897844203874
$ws.Range("P8").Value = 0.8730897844203874
$ws.Range("Q8").Value = 609.9411813998399
$ws.Range("R8").Value = 5489.470632598559
$ws.Range("S8").Value = 0.7544179441719365
$ws.Range("T8").Value = 0.7544179441719365

$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Il34"
$ws.Range("C9").Value = "Csf1r"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 27.63137633333334
$ws.Range("H9").Value = 82.894129
$ws.Range("I9").Value = 0.8640783085931617
$ws.Range("J9").Value = 0.8640783085931617
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1.456833333333333
$ws.Range("N9").Value = 4.3705
$ws.Range("O9").Value = 0.05762133790721463
$ws.Range("P9").Value = 0.05762133790721465
$ws.Range("Q9").Value = 40.25431008827778
$ws.Range("R9").Value = 362.2887907945
$ws.Range("S9").Value = 0.04978934819774105
$ws.Range("T9").Value = 0.04978934819774106

$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Il34"
$ws.Range("C10").Value = "Csf1r"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 27.63137633333334
$ws.Range("H10").Value = 82.894129
$ws.Range("I10").Value = 0.8640783085931617
$ws.Range("J10").Value = 0.8640783085931617
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 1.751822333333333
$ws.Range("N10").Value = 5.255467
$ws.Range("O10").Value = 0.06928887767239802
$ws.Range("P10").Value = 0.06928887767239804
$ws.Range("Q10").Value = 48.40526216147145
$ws.Range("R10").Value = 435.647359453243
$ws.Range("S10").Value = 0.05987101622348417
$ws.Range("T10").Value = 0.05987101622348418
